# Add two new columns ("NIF" and "pollingStation") to the citizens table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("H1").Value = "NIF"
$ws.Range("I1").Value = "pollingStation"

# Data rows: pollingStation (and the demo NIF placeholder) are 1, 2, 3 for
# the three citizens already on the sheet.
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1

$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 2

$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3

# Leave the selection on the last written cell, like a user would after
# typing in the new column.
$ws.Range("G4").Select()
